$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "is_global_admin" header (column AH) and its data column are removed
# from the "user" sheet. Deleting the entire column shifts every column to
# its right one place to the left, which also drops the now-unused
# "is_global_admin" shared string and renumbers the shared-string indices
# used by the remaining cells - matching the target workbook exactly.
$ws.Columns("AH").Delete()

# Reflect where the user's selection ended up after performing the deletion.
$null = $ws.Range("AK2").Select()
